$wb = $excel.ActiveWorkbook
$wb.Names.Add("test4", "=Resúmen!`$W15:`$W18")
$wb.Names.Add("test5", "=Resúmen!W`$15:W`$18")
$wb.Names.Add("test6", "=Resúmen!`$W`$15:`$W`$18")
